# Update the one quiz result that changed between the 05:26:42 and
# 05:40:08 submissions for "Deepa" (row 3 of the "Quiz Results" sheet):
#   - I-Number goes from I079692 -> I079693
#   - Score goes from 60 -> 100
#   - Correct Answers goes from 3 -> 5
#   - Duration (sec) goes from 35 -> 33
#   - Completion Date goes from 2025-04-25 05:26:42 -> 2025-04-25 05:40:08
# Total Questions (5) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quiz Results")

$ws.Range("C3").Value = "I079693"
$ws.Range("D3").Value = 100.0
$ws.Range("F3").Value = 5.0
$ws.Range("G3").Value = 33.0
$ws.Range("H3").Value = "2025-04-25 05:40:08"
